$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells AD1:AF1 with same style as existing header (copy formatting from A1)
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Fill season record (Wins, Losses, Ties) for every data row (2 through 44)
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 72
    $ws.Cells.Item($r, 31).Value2 = 90
    $ws.Cells.Item($r, 32).Value2 = 0
}
